$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.412.18"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").Value = "1.654.54"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.50"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.05"
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "1.888.98"
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("D13").Value = "1.663.15"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.575"
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.08"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.82"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "27.426.96"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.05"
$ws.Range("E18").Value = "  -7.94%  "
$ws.Range("E19").Value = "  -2.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.51"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.38"
$ws.Range("E22").Value = "  -3.42%  "
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("E24").Value = "  -2.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.94"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.21"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("E27").Value = "  -3.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0497"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  -4.39%  "
$ws.Range("E32").Value = "  -2.57%  "
$ws.Range("D33").Value = "1.466.37"
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("E35").Value = "  -4.22%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("E38").Value = "  -3.43%  "
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.05"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.46"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.39"
$ws.Range("E43").Value = "  -6.54%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "1.797.16"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.781"
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.38"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").Value = "0.0₆0106"
$ws.Range("E49").Value = "  -4.59%  "
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.78"
$ws.Range("E51").Value = "  -0.82%  "
